$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content beyond new range, and full sheet first to avoid stale data
$ws.Cells.Clear()

# Write all cell values
$ws.Range('A1').Value = 'Test: Standardablauf/Epic'
$ws.Range('A2').Value = 'Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 1.'
$ws.Range('B2').Value = 'Es wird im DropDown-Menü eine Liste aller möglichen Spieler angezeigt.'
$ws.Range('A3').Value = 'Der Nutzer wählt die Option "Mensch" im DropDown-Menü aus.'
$ws.Range('B3').Value = 'Als Spieler 1 wird ein Mensch festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde.'
$ws.Range('A4').Value = 'Der Nutzer klickt auf das DropDown-Menü für die Auswahl von Spieler 2.'
$ws.Range('B4').Value = 'Es wird im DropDown-Menü eine Liste aller möglichen Spieler angezeigt.'
$ws.Range('A5').Value = 'Der Nutzer wählt die Option "KI 1" im DropDown-Menü aus.'
$ws.Range('B5').Value = 'Als Spieler 1 wird die "KI 1" festgelegt. Im DropDown-Menü wird angezeigt, dass die Option ausgewählt wurde.'
$ws.Range('A6').Value = 'Der Nutzer klickt auf das Feld "Spiel starten".'
$ws.Range('B6').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A7').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B7').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A8').Value = 'Der Nutzer klickt auf das "NextMove"-Symbol (="Play"-Symbol mit | rechts).'
$ws.Range('B8').Value = 'Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A9').Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range('B9').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'
$ws.Range('A10').Value = 'Der Nutzer klickt auf ein zufälliges freies Spielfeld.'
$ws.Range('B10').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt. Ein zufälliges freies Feld wird von der KI besetzt und mit Zeichen "O" versehen.  Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A11').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf den Zug der KI zu warten, bis das Spiel vorbei ist.'
$ws.Range('B11').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt. Es wird in die Belohnungsansicht gewechselt. Auf der rechten Seite wird als Graph der gesamte gewichtete Verlauf mit pro Zustand allen möglichen Äquivalenzklassenvertretern der Folgezuständen angezeigt.'
$ws.Range('A12').Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range('B12').Value = 'Es wird in den Startansicht gewechselt, in welchem die Spieler ausgewählt werden.'
$ws.Range('A13').Value = 'Der Nutzer klickt auf die Checkbox "Startansicht überspringen".'
$ws.Range('B13').Value = 'Die Checkbox wird mit einem Haken versehen.'
$ws.Range('A14').Value = 'Der Nutzer klickt auf das Feld "Spiel starten".'
$ws.Range('B14').Value = 'Es wird in die Spielansicht gewechselt. Das "Play"-Symbol ist weiterhin durch ein "Pause"-Symbol ersetzt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Äquivalenzklassenvertreter der Folgezustände an.'
$ws.Range('A15').Value = ' Der Nutzer klickt auf die Checkbox "Belohnungsansicht überspringen".'
$ws.Range('B15').Value = 'Die Checkbox wird mit einem Haken versehen.'
$ws.Range('A16').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf den Zug der KI zu warten, bis das Spiel vorbei ist.'
$ws.Range('B16').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird direkt das Spielfeld geleert und der Verlauf-Graph auf das leere Feld mit allen möglichen Äquivalenzklassenvertretern der Folgezuständen geändert.'
$ws.Range('A17').Value = 'Der Nutzer klickt auf die Checkbox "Startansicht überspringen".'
$ws.Range('B17').Value = 'Der Haken in der Checkbox wird entfernt.'
$ws.Range('A18').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf den Zug der KI zu warten, bis das Spiel vorbei ist.'
$ws.Range('B18').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird in die Startansicht gewechselt.'
$ws.Range('A19').Value = 'Der Nutzer klickt auf das Feld "Reset" neben KI 1.'
$ws.Range('B19').Value = 'Es ändert sich nichts, da nur die Gewichte der KI gelöscht wurden.'
$ws.Range('A21').Value = 'Test: AutoPlay/Pause'
$ws.Range('A22').Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "KI 1" und klickt auf das Feld "Spiel starten".'
$ws.Range('B22').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A23').Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range('B23').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'
$ws.Range('A24').Value = 'Der Nutzer klickt auf das "NextMove"-Symbol.'
$ws.Range('B24').Value = 'Das "Pause"-Symbol wird durch ein "Play"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'
$ws.Range('A25').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten rechts.'
$ws.Range('B25').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A26').Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range('B26').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt. Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A27').Value = 'Der Nutzer klickt auf das "NextMove"-Symbol.'
$ws.Range('B27').Value = 'Das "Pause"-Symbol wird durch ein "Play"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'
$ws.Range('A28').Value = 'Der Nutzer klickt auf ein zufälliges unbesetztes Feld des Spielfeldes.'
$ws.Range('B28').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A29').Value = 'Der Nutzer klickt auf das "NextMove"-Symbol.'
$ws.Range('B29').Value = 'Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A30').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf das "NextMove"-Symbol zu klicken, bis das Spiel vorbei ist und Spieler 1 gewonnen hat.'
$ws.Range('B30').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt. Es wird in die Belohnungsansicht gewechselt. Auf der rechten Seite wird als Graph der gesamte gewichtete Verlauf mit pro Zustand allen möglichen Äquivalenzklassenvertretern der Folgezuständen angezeigt.'
$ws.Range('A31').Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range('B31').Value = 'Es wird in den Startansicht gewechselt, in welchem die Spieler ausgewählt werden.'
$ws.Range('A32').Value = 'Der Nutzer klickt auf das Feld "Spiel starten".'
$ws.Range('B32').Value = 'Es wird in die Spielansicht gewechselt. Das "Play"-Symbol ist zu sehen und nicht durch das "Pause"-Symbol ersetzt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Äquivalenzklassenvertreter der Folgezustände an.'
$ws.Range('A33').Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range('B33').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt. Es wird kein Zug ausgeführt, da der Mensch am Zug ist.'
$ws.Range('A34').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld rechts oben.'
$ws.Range('B34').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt. Ein zufälliges freies Feld wird von der KI besetzt und mit Zeichen "O" versehen.  Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A35').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld rechts oben.'
$ws.Range('B35').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt. Ein zufälliges freies Feld wird von der KI besetzt und mit Zeichen "O" versehen.  Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A36').Value = 'Der Nutzer klickt auf das "Pause"-Symbol.'
$ws.Range('B36').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt.'
$ws.Range('A37').Value = 'Der Nutzer klickt auf ein freies Feld des Spielfeldes, welches das Spiel nicht beendet. (Ein Zug beendet ein Spiel, falls drei Kreuze in in einer waagrechten/senkrechten/schrägen Reihe sind)'
$ws.Range('B37').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt. Die KI macht keinen Zug, da AutoPlay mittels "Pause"-Button deaktiviert wurde.'
$ws.Range('A38').Value = 'Der Nutzer klickt auf das "Play"-Symbol.'
$ws.Range('B38').Value = 'Das "Play"-Symbol wird durch ein "Pause"-Symbol ersetzt. Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A39').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf den Zug der KI zu warten, bis das Spiel vorbei ist.'
$ws.Range('B39').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt. Es wird in die Belohnungsansicht gewechselt. Auf der rechten Seite wird als Graph der gesamte gewichtete Verlauf mit pro Zustand allen möglichen Äquivalenzklassenvertretern der Folgezuständen angezeigt.'
$ws.Range('A40').Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range('B40').Value = 'Es wird in die Startansicht gewechselt, in welchem die Spieler ausgewählt werden.'
$ws.Range('A41').Value = 'Der Nutzer klickt auf das Feld "Spiel starten".'
$ws.Range('B41').Value = 'Es wird in die Spielansicht gewechselt. Das "Play"-Symbol ist weiterhin durch ein "Pause"-Symbol ersetzt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle möglichen Äquivalenzklassenvertreter der Folgezustände an.'
$ws.Range('A43').Value = 'Test: Spiel mit zwei Menschen'
$ws.Range('A44').Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "Mensch" und klickt auf das Feld "Spiel starten".'
$ws.Range('B44').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A45').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B45').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A46').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B46').Value = 'Es passiert nichts, da das Feld bereits besetzt ist.'
$ws.Range('A47').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten links.'
$ws.Range('B47').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A48').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben in der Mitte.'
$ws.Range('B48').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A49').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten in der Mitte.'
$ws.Range('B49').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A50').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben rechts.'
$ws.Range('B50').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf ohne Äquivalenzklassenvertreter der Folgezustände angezeigt. Das Spielergebnis wird angezeigt: "Spieler 1 gewinnt!"'
$ws.Range('A52').Value = 'Test: Spielfeldbesetzung - Randfälle'
$ws.Range('A53').Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "KI 1" und klickt auf das Feld "Spiel starten".'
$ws.Range('B53').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A54').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B54').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A55').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben rechts.'
$ws.Range('B55').Value = 'Es passiert nichts, da der Spieler "Mensch" nicht an der Reihe ist.'
$ws.Range('A56').Value = 'Der Nutzer klickt auf das "NextMove"-Symbol.'
$ws.Range('B56').Value = 'Die KI macht einen zufälligen Zug und belegt somit ein zufälliges unbesetztes Feld. Dieses wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A57').Value = 'Der Nutzer klickt auf dem Spielfeld auf das gerade von der KI ausgewählte Feld, welches mit einem "O" versehen wurde.'
$ws.Range('B57').Value = 'Es passiert nichts, da das Feld bereits besetzt ist.'
$ws.Range('A58').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B58').Value = 'Es passiert nichts, da das Feld bereits besetzt ist.'
$ws.Range('A59').Value = 'Der Nutzer wechselt damit ab, auf ein zufälliges freies Spielfeld zu klicken und auf das "NextMove"-Symbol zu klicken, bis das Spiel vorbei ist und Spieler 1 gewonnen hat.'
$ws.Range('B59').Value = 'Das Spielfeld wird laufend aktualisiert. Der Graph zeigt jeweils immer den aktualisierten Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände an. Sobald das Spiel vorbei ist, wird das Spielergebnis angezeigt. Es wird in die Belohnungsansicht gewechselt. Auf der rechten Seite wird als Graph der gesamte gewichtete Verlauf mit pro Zustand allen möglichen Äquivalenzklassenvertretern der Folgezuständen angezeigt.'
$ws.Range('A60').Value = 'Der Nutzer klickt auf das Feld "Belohnung ausführen".'
$ws.Range('B60').Value = 'Es wird in den Startansicht gewechselt, in welchem die Spieler ausgewählt werden.'
$ws.Range('A61').Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "Mensch" und klickt auf das Feld "Spiel starten".'
$ws.Range('B61').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A62').Value = 'Führe Schritte von "Test: Spiel mit zwei Menschen" durch.'
$ws.Range('A63').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten rechts.'
$ws.Range('B63').Value = 'Es passiert nichts, da das Spiel bereits beendet ist.'
$ws.Range('A65').Value = 'Test: Unentschieden'
$ws.Range('A66').Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "Mensch" und klickt auf das Feld "Spiel starten".'
$ws.Range('B66').Value = 'Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an.'
$ws.Range('A67').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben links.'
$ws.Range('B67').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A68').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld in der Mitte links.'
$ws.Range('B68').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A69').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben in der Mitte.'
$ws.Range('B69').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A70').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld in der Mitte.'
$ws.Range('B70').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A71').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld in der Mitte rechts.'
$ws.Range('B71').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A72').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld oben rechts.'
$ws.Range('B72').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A73').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten links.'
$ws.Range('B73').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A74').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten rechts.'
$ws.Range('B74').Value = 'Das angeklickte Feld wird mit dem Zeichen "O" versehen. Im Graph wird der aktualisierte Verlauf mit allen für die aktuelle Spielsituation möglichen Äquivalenzklassenvertretern der Folgezustände angezeigt.'
$ws.Range('A75').Value = 'Der Nutzer klickt auf dem Spielfeld auf das Feld unten in der Mitte.'
$ws.Range('B75').Value = 'Das angeklickte Feld wird mit dem Zeichen "X" versehen. Im Graph wird der aktualisierte Verlauf ohne Äquivalenzklassenvertreter der Folgezustände angezeigt. Das Spielergebnis wird angezeigt: "Unentschieden!"'

# Bold section header rows
$ws.Range('A1').Font.Bold = $true
$ws.Range('A21').Font.Bold = $true
$ws.Range('A43').Font.Bold = $true
$ws.Range('A52').Font.Bold = $true
$ws.Range('A65').Font.Bold = $true

# Update selection to match target view state
$ws.Range('B11').Select()
